$d = $word.ActiveDocument

# Locate the scope-statement paragraph that ends with the "(232 words)" count.
# (There is a near-duplicate paragraph later in the document that must NOT be touched.)
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*(232 words)*") {
        $target = $p
    }
}

# 1) Remove ", with implications for the fate of carbon exported to depth" (scoped to this paragraph only)
$rng1 = $target.Range
$rng1.Find.Execute(", with implications for the fate of carbon exported to depth.", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2)

# 2) Remove " and place an emphasis on reproducibility by collecting an entirely separate set of data from a new location a year later for comparison" (scoped)
$rng2 = $target.Range
$rng2.Find.Execute(" and place an emphasis on reproducibility by collecting an entirely separate set of data from a new location a year later for comparison.", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2)

# 3) Update the word count from 232 to 199, keeping it as its own run (split off from
#    surrounding text), matching the structure produced by the original edit.
$numRange = $target.Range.Duplicate
$numRange.Find.Execute("232")
$numRange.Text = "199"
$numRange.Bold = 1
$numRange.Bold = 0

Write-Output $target.Range.Text
